$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the leave balance start/end date text values (stored as text, not real dates)
$ws.Range("A2").Value = "2023-August-1"
$ws.Range("B2").Value = "2023-August-31"

# Update the active cell selection shown on the sheet view
$ws.Range("D15").Select()
